# "Fixing bug with CODE.SUBST."
# Adds a new "Subst" worksheet as the first sheet in the workbook, populates
# it with a small instruction/closing/blocks scratch table used to debug the
# CODE.SUBST instruction, and makes it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# --- create the new sheet and move it to the front of the tab strip ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Subst"
$newSheet.Move($wb.Worksheets.Item(1))

$ws = $wb.Worksheets.Item(1)

# --- header row ---
$ws.Range("A1").Value = "Instruction"
$ws.Range("B1").Value = "closing"
$ws.Range("C1").Value = "extra_blocks"
$ws.Range("D1").Value = "wanted_blocks"
$ws.Range("E1").Value = "wanted_stack"
$ws.Range("F1").Value = "item_number"

# --- data rows ---
$ws.Range("A2").Value = "{:instruction EXEC.NOOP_OPEN_PAREN :close 0}\"
$ws.Range("B2").Value = -1
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = "{:instruction EXEC.DO*RANGE :close 0}\"
$ws.Range("B3").Value = -1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "0, 1"
$ws.Range("F3").Value = 0

$ws.Range("A4").Value = "{:instruction FLOAT.+ :close 2}\"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0

$ws.Range("A6").Value = "{:instruction EXEC.NOOP_OPEN_PAREN :close 0}\"
$ws.Range("A7").Value = "{:instruction FLOAT.+ :close 2}\"
$ws.Range("A8").Value = "{:instruction FLOAT.- :close 1}\"
$ws.Range("A9").Value = "{:instruction FLOAT.+ :close 1}\"
$ws.Range("A10").Value = "{:instruction CODE.SUBST :close 0}\"

# --- column E is right aligned in the original workbook's similar tables ---
$ws.Range("E1:E4").HorizontalAlignment = -4152

# --- column widths (best-fit look) ---
$ws.Columns.Item(1).ColumnWidth = 42
$ws.Columns.Item(2).ColumnWidth = 6.625
$ws.Columns.Item(3).ColumnWidth = 11.125
$ws.Columns.Item(4).ColumnWidth = 13.125
$ws.Columns.Item(5).ColumnWidth = 12.125
$ws.Columns.Item(6).ColumnWidth = 11.875

# --- make "Subst" the active sheet/tab, with A9 selected, zoomed in ---
$ws.Select()
[void]($excel.ActiveWindow.Zoom = 180)
[void]($ws.Range("A9").Select())
